# Weekly price-sheet update: insert one new "Perejil" (parsley) price
# observation as row 271 of the only worksheet, pushing every existing
# row from 271 down to 393 one row further (to 272..394).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 271; rows 271..393 shift down to 272..394.
$ws.Rows.Item(271).Insert()

# Populate the newly inserted row 271 with the new weekly observation.
$ws.Cells.Item(271, 1).Value  = 10
$ws.Cells.Item(271, 2).Value  = 'Vega Modelo de Temuco'
$ws.Cells.Item(271, 3).Value  = 'La Araucanía'
$ws.Cells.Item(271, 4).Value  = 44839
$ws.Cells.Item(271, 5).Value  = 9
$ws.Cells.Item(271, 6).Value  = 100112044
$ws.Cells.Item(271, 7).Value  = 'Perejil'
$ws.Cells.Item(271, 8).Value  = 'Sin especificar'
$ws.Cells.Item(271, 9).Value  = 'Primera'
$ws.Cells.Item(271, 10).Value = 80
$ws.Cells.Item(271, 11).Value = 3300
$ws.Cells.Item(271, 12).Value = 3300
$ws.Cells.Item(271, 13).Value = 3300
$ws.Cells.Item(271, 14).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(271, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(271, 16).Value = 1100
$ws.Cells.Item(271, 17).Value = 3
$ws.Cells.Item(271, 18).Value = 'Hortaliza'
